# The deck's design was switched from the "Integral" (Red Violet) theme to
# the built-in "Office Theme" (Office palette). In the OOXML this shows up
# as the slide master's theme part (ppt/theme/theme1.xml) taking on the
# Office Theme's colour scheme (the font scheme and format scheme were
# already identical between the two themes, so only the 12 theme colours
# actually change).
$p = $ppt.ActivePresentation

# RGB() isn't available in this host, so build the BGR-packed COM colour
# value (0x00BBGGRR) from R/G/B components ourselves.
function ComRGB([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

$officeTheme = @{
    1  = (ComRGB 0x00 0x00 0x00)   # dk1
    2  = (ComRGB 0xFF 0xFF 0xFF)   # lt1
    3  = (ComRGB 0x44 0x54 0x6A)   # dk2
    4  = (ComRGB 0xE7 0xE6 0xE6)   # lt2
    5  = (ComRGB 0x5B 0x9B 0xD5)   # accent1
    6  = (ComRGB 0xED 0x7D 0x31)   # accent2
    7  = (ComRGB 0xA5 0xA5 0xA5)   # accent3
    8  = (ComRGB 0xFF 0xC0 0x00)   # accent4
    9  = (ComRGB 0x44 0x72 0xC4)   # accent5
    10 = (ComRGB 0x70 0xAD 0x47)   # accent6
    11 = (ComRGB 0x05 0x63 0xC1)   # hlink
    12 = (ComRGB 0x95 0x4F 0x72)   # folHlink
}

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeTheme[$i]
}
